$d = $word.ActiveDocument

# Find the paragraph containing the sentence about the ttest between
# Door3 and Switch.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Unlike the ttest between Door3 and DoesntMatter*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $paraRange = $target.Range

    # Remove " or 0.99" so the sentence reads "...consistently 1.0 over
    # multiple tests." instead of "...consistently 1.0 or 0.99 over
    # multiple tests."
    $fullText = $paraRange.Text
    $idx = $fullText.IndexOf(" or 0.99")
    if ($idx -ge 0) {
        $delStart = $paraRange.Start + $idx
        $delEnd = $delStart + 8
        $d.Range($delStart, $delEnd).Text = ""
    }

    # The author's original text was split into two runs at the boundary
    # right before "over multiple tests."; recreate that run break by
    # nudging character formatting on/off at that same point (this does
    # not change the visible formatting, just forces the run split).
    $newText = $target.Range.Text
    $splitIdx = $newText.IndexOf("over multiple tests.")
    if ($splitIdx -ge 0) {
        $splitStart = $target.Range.Start + $splitIdx
        $splitEnd = $target.Range.End
        $secondRun = $d.Range($splitStart, $splitEnd)
        $secondRun.Font.Bold = 1
        $secondRun.Font.Bold = 0
    }
}
